$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FLO_SUBSIDY")

# Update the NL limit value for FLO_SUB / Year 2050 / ELC_FIN_DEM / ELC_GRID_RES from 2 to 1
$ws.Range("F4").Value = 1

# Leave the cursor where the author left it before saving
[void]$ws.Range("B12").Select()
